# edit.ps1
# Applies: add 35 new digimon rows, correct one duplicate text,
# sort the B:C table by offset (col B) ascending, shrink the
# hidden _FilterDatabase name back to the header row, and leave
# the selection where the author's cursor ended up (B53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the pre-existing duplicate-name row (offset 2611 was
#        mistakenly sharing "타네몬" with offset 2511; give it its
#        own name) --------------------------------------------------
$ws.Range("C18").Value = "푸카몬"

# --- 2) Append the new offset/name pairs below the existing table -
$newData = @(
  @(2131, "그레이몬"),
  @(2141, "메탈그레이몬"),
  @(2151, "워그레이몬"),
  @(2221, "가부몬"),
  @(2231, "가루루몬"),
  @(2241, "워가루루몬"),
  @(2251, "메탈가루루몬"),
  @(2321, "피요몬"),
  @(2331, "버드라몬"),
  @(2341, "가루다몬"),
  @(2351, "호우오우몬"),
  @(2431, "캅테리몬"),
  @(2441, "아트라캅테리몬"),
  @(2451, "헤라클캅테리몬"),
  @(2531, "토게몬"),
  @(2541, "리리몬"),
  @(2551, "로제몬"),
  @(2621, "고마몬"),
  @(2631, "잇카쿠몬"),
  @(2641, "즈도몬"),
  @(2651, "바이크몬"),
  @(2731, "엔제몬"),
  @(2741, "홀리엔제몬"),
  @(2751, "세라피몬"),
  @(2821, "플롯트몬"),
  @(2831, "테일몬"),
  @(2841, "엔제우몬"),
  @(2851, "오파니몬"),
  @(3021, "에테몬"),
  @(3031, "반데몬"),
  @(4182, "가지몬"),
  @(4191, "코카토리몬"),
  @(4201, "핏코로몬"),
  @(4905, "티라노몬"),
  @(5011, "겐나이")
)

$r = 32
foreach ($item in $newData) {
  $ws.Cells.Item($r, 2).Value = $item[0]
  $ws.Cells.Item($r, 3).Value = $item[1]
  $r = $r + 1
}

# New text cells need the same "@" text style as the rest of column C
$ws.Range("C32:C66").NumberFormat = "@"

# --- 3) Sort the whole table (still headerless row 1) by offset ---
$ws.Range("B2:C66").Sort($ws.Range("B2"), 1, $null, $null, 1, $null, $null, 2)

# --- 4) The filter-database bookmark collapses back to just the
#        header row once the filter is cleared/reset -------------
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$B`$1:`$C`$1"
  }
}

# --- 5) Leave the cursor where the author left it -----------------
$ws.Range("B53").Select()
